$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (row 1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 20 de Agosto de 2020 a las 14:30"

# Reorder "Montserrat" / "Islas Malvinas" rows (Islas Malvinas now comes first)
$ws.Cells.Item(213,1).Value = "Islas Malvinas"
$ws.Cells.Item(214,1).Value = "Montserrat"

# Row 4 - Estados Unidos
$ws.Cells.Item(4,2).Value = 5701878
$ws.Cells.Item(4,3).Value = 947
$ws.Cells.Item(4,4).Value = 3063252
$ws.Cells.Item(4,5).Value = 2462249
$ws.Cells.Item(4,7).Value = 40
$ws.Cells.Item(4,8).Value = 176377

# Row 22 - Alemania
$ws.Cells.Item(22,2).Value = 229826
$ws.Cells.Item(22,3).Value = 126
$ws.Cells.Item(22,5).Value = 16611
$ws.Cells.Item(22,7).Value = 1
$ws.Cells.Item(22,8).Value = 9315

# Row 36 - Suecia
$ws.Cells.Item(36,2).Value = 85810
$ws.Cells.Item(36,7).Value = 5
$ws.Cells.Item(36,8).Value = 5805

# Row 41 - Kuwait
$ws.Cells.Item(41,2).Value = 78767
$ws.Cells.Item(41,3).Value = 622
$ws.Cells.Item(41,4).Value = 70642
$ws.Cells.Item(41,5).Value = 7616
$ws.Cells.Item(41,7).Value = 2
$ws.Cells.Item(41,8).Value = 509

# Row 44 - Emiratos Arabes Unidos
$ws.Cells.Item(44,2).Value = 65802
$ws.Cells.Item(44,3).Value = 461
$ws.Cells.Item(44,4).Value = 58153
$ws.Cells.Item(44,5).Value = 7280
$ws.Cells.Item(44,7).Value = 2
$ws.Cells.Item(44,8).Value = 369

# Row 61 - Uzbekistan
$ws.Cells.Item(61,2).Value = 37547
$ws.Cells.Item(61,3).Value = 435
$ws.Cells.Item(61,4).Value = 33261
$ws.Cells.Item(61,5).Value = 4034
$ws.Cells.Item(61,7).Value = 4
$ws.Cells.Item(61,8).Value = 252

# Row 78 - Bosnia y Herzegovina
$ws.Cells.Item(78,2).Value = 17029
$ws.Cells.Item(78,3).Value = 338
$ws.Cells.Item(78,4).Value = 10881
$ws.Cells.Item(78,5).Value = 5633
$ws.Cells.Item(78,7).Value = 8
$ws.Cells.Item(78,8).Value = 515

# Row 80 - Dinamarca
$ws.Cells.Item(80,2).Value = 16056
$ws.Cells.Item(80,3).Value = 116
$ws.Cells.Item(80,4).Value = 13769
$ws.Cells.Item(80,5).Value = 1666

# Row 83 - Republica de Macedonia
$ws.Cells.Item(83,2).Value = 13194
$ws.Cells.Item(83,3).Value = 118
$ws.Cells.Item(83,4).Value = 9752
$ws.Cells.Item(83,5).Value = 2888
$ws.Cells.Item(83,7).Value = 3
$ws.Cells.Item(83,8).Value = 554

# Row 89 - Noruega
$ws.Cells.Item(89,5).Value = 1041
$ws.Cells.Item(89,7).Value = 2
$ws.Cells.Item(89,8).Value = 264

# Row 102 - Croacia
$ws.Cells.Item(102,2).Value = 7329
$ws.Cells.Item(102,3).Value = 255
$ws.Cells.Item(102,4).Value = 5472
$ws.Cells.Item(102,5).Value = 1689

# Row 136 - Islandia
$ws.Cells.Item(136,2).Value = 2040
$ws.Cells.Item(136,3).Value = 5
$ws.Cells.Item(136,4).Value = 1913
$ws.Cells.Item(136,5).Value = 117

# Row 158 - Vietnam
$ws.Cells.Item(158,2).Value = 1007
$ws.Cells.Item(158,3).Value = 13
$ws.Cells.Item(158,4).Value = 542
$ws.Cells.Item(158,5).Value = 440

# Row 194 - Liechtenstein
$ws.Cells.Item(194,2).Value = 99
$ws.Cells.Item(194,3).Value = 1
$ws.Cells.Item(194,5).Value = 10

# Row 213 - Islas Malvinas (was Montserrat's position)
$ws.Cells.Item(213,4).Value = 13
$ws.Cells.Item(213,8).Value = 0

# Row 214 - Montserrat (was Islas Malvinas's position)
$ws.Cells.Item(214,4).Value = 12
$ws.Cells.Item(214,8).Value = 1
